$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F3 (12 -> 13) and F4 (965 -> 966)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13
$ws1.Range("F4").Value = 966

# Sheet "全部类型" - update F3 (12 -> 13) and F4 (965 -> 966)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 13
$ws4.Range("F4").Value = 966
